$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1754
$ws.Range("I2").Value = 1754
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1754
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1641
$ws.Range("H9").Value = 109.55556
$ws.Range("I9").Value = 109.55556
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 109.55556
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = 59.44444
$ws.Range("H64").Value = 5031.25
$ws.Range("J64").Value = 4985.7144
$ws.Range("L64").Value = 4985.7144
$ws.Range("N64").Value = -5481.7144
$ws.Range("H67").Value = 5031.25
$ws.Range("J67").Value = 4985.7144
$ws.Range("L67").Value = 4985.7144
$ws.Range("N67").Value = -6701.7144
$ws.Range("H100").Value = 1719.8
$ws.Range("I100").Value = 1649.75
$ws.Range("K100").Value = 1649.75
$ws.Range("M100").Value = -1108.75
$ws.Range("H113").Value = 4250
$ws.Range("I113").Value = 2875
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 2875
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = 379
$ws.Range("N113").Value = -13508
$ws.Range("H135").Value = 17354.5
$ws.Range("I135").Value = 825.4
$ws.Range("K135").Value = 7428.599999999999
$ws.Range("M135").Value = -4893.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1000.6429
$ws.Range("I2").Value = 1100.8334
$ws.Range("J2").Value = 399.5
$ws.Range("K2").Value = 1100.8334
$ws.Range("L2").Value = 399.5
$ws.Range("M2").Value = -987.8334
$ws.Range("N2").Value = -625.5
$ws.Range("H74").Value = 27800928
$ws.Range("J74").Value = 26044
$ws.Range("L74").Value = 26044
$ws.Range("N74").Value = -27792
$ws.Range("H77").Value = 27800928
$ws.Range("J77").Value = 26044
$ws.Range("L77").Value = 130220
$ws.Range("N77").Value = -138956
$ws.Range("H102").Value = 16421.166
$ws.Range("I102").Value = 19478.732
$ws.Range("J102").Value = 1133.3334
$ws.Range("K102").Value = 19478.732
$ws.Range("L102").Value = 1133.3334
$ws.Range("M102").Value = -17856.732
$ws.Range("N102").Value = -4377.3334
$ws.Range("H103").Value = 48996.668
$ws.Range("J103").Value = 48996.668
$ws.Range("L103").Value = 48996.668
$ws.Range("N103").Value = -51340.668
$ws.Range("H116").Value = 1000.6429
$ws.Range("I116").Value = 1100.8334
$ws.Range("J116").Value = 399.5
$ws.Range("K116").Value = 1100.8334
$ws.Range("L116").Value = 399.5
$ws.Range("M116").Value = 1193.1666
$ws.Range("N116").Value = -4987.5
$ws.Range("H121").Value = 109985
$ws.Range("J121").Value = 109985
$ws.Range("L121").Value = 109985
$ws.Range("N121").Value = -113479
$ws.Range("H132").Value = 12178.467
$ws.Range("I132").Value = 5713.857
$ws.Range("J132").Value = 17835
$ws.Range("K132").Value = 17141.571
$ws.Range("L132").Value = 53505
$ws.Range("M132").Value = -14611.571
$ws.Range("N132").Value = -58565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1000.6429
$ws.Range("I3").Value = 1100.8334
$ws.Range("J3").Value = 399.5
$ws.Range("K3").Value = 1100.8334
$ws.Range("L3").Value = 399.5
$ws.Range("M3").Value = -986.8334
$ws.Range("N3").Value = -627.5
$ws.Range("H86").Value = 2560
$ws.Range("I86").Value = 2560
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2560
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -1437
$ws.Range("H89").Value = 2560
$ws.Range("I89").Value = 2560
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12800
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -7184
$ws.Range("H99").Value = 2090.1875
$ws.Range("I99").Value = 1688.7858
$ws.Range("K99").Value = 1688.7858
$ws.Range("M99").Value = -190.7858000000001
$ws.Range("H105").Value = 3462.375
$ws.Range("I105").Value = 1849.5
$ws.Range("K105").Value = 1849.5
$ws.Range("M105").Value = -102.5
$ws.Range("H107").Value = 1557.4
$ws.Range("I107").Value = 1557.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1557.4
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 362.5999999999999
$ws.Range("H134").Value = 180532.83
$ws.Range("I134").Value = 6584.5
$ws.Range("J134").Value = 267507
$ws.Range("K134").Value = 19753.5
$ws.Range("L134").Value = 802521
$ws.Range("M134").Value = -17218.5
$ws.Range("N134").Value = -807591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1160.0834
$ws.Range("I58").Value = 876
$ws.Range("J58").Value = 2012.3334
$ws.Range("K58").Value = 876
$ws.Range("L58").Value = 2012.3334
$ws.Range("M58").Value = -673
$ws.Range("N58").Value = -2418.3334
$ws.Range("H99").Value = 2753.5386
$ws.Range("I99").Value = 2531.8572
$ws.Range("J99").Value = 3012.1667
$ws.Range("K99").Value = 2531.8572
$ws.Range("L99").Value = 3012.1667
$ws.Range("M99").Value = -1033.8572
$ws.Range("N99").Value = -6008.1667
$ws.Range("H126").Value = 2753.5386
$ws.Range("I126").Value = 2531.8572
$ws.Range("J126").Value = 3012.1667
$ws.Range("K126").Value = 7595.571599999999
$ws.Range("L126").Value = 9036.500100000001
$ws.Range("M126").Value = -5125.571599999999
$ws.Range("N126").Value = -13976.5001
$ws.Range("H132").Value = 2883.2222
$ws.Range("I132").Value = 2618.625
$ws.Range("K132").Value = 7855.875
$ws.Range("M132").Value = -5325.875
$ws.Range("H135").Value = 75333.336
$ws.Range("I135").Value = 40000
$ws.Range("J135").Value = 93000
$ws.Range("K135").Value = 40000
$ws.Range("L135").Value = 93000
$ws.Range("M135").Value = -34930
$ws.Range("N135").Value = -103140
$ws.Range("H136").Value = 1160.0834
$ws.Range("I136").Value = 876
$ws.Range("J136").Value = 2012.3334
$ws.Range("K136").Value = 2628
$ws.Range("L136").Value = 6037.0002
$ws.Range("M136").Value = -78
$ws.Range("N136").Value = -11137.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5322.091
$ws.Range("I102").Value = 3282.5557
$ws.Range("K102").Value = 3282.5557
$ws.Range("M102").Value = -1660.5557
$ws.Range("H113").Value = 4266.1113
$ws.Range("I113").Value = 4071.4285
$ws.Range("K113").Value = 4071.4285
$ws.Range("M113").Value = -1901.4285
$ws.Range("H132").Value = 58840644
$ws.Range("J132").Value = 52500
$ws.Range("L132").Value = 157500
$ws.Range("N132").Value = -162560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5724.9473
$ws.Range("I122").Value = 4555.2856
$ws.Range("K122").Value = 13665.8568
$ws.Range("M122").Value = -11215.8568
$ws.Range("H132").Value = 97846.766
$ws.Range("I132").Value = 11950.3
$ws.Range("J132").Value = 384168.34
$ws.Range("K132").Value = 35850.89999999999
$ws.Range("L132").Value = 1152505.02
$ws.Range("M132").Value = -33320.89999999999
$ws.Range("N132").Value = -1157565.02
$ws.Range("H136").Value = 35001.895
$ws.Range("I136").Value = 4924.6665
$ws.Range("J136").Value = 86562.86
$ws.Range("K136").Value = 14773.9995
$ws.Range("L136").Value = 259688.58
$ws.Range("M136").Value = -12223.9995
$ws.Range("N136").Value = -264788.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9762.6
$ws.Range("I122").Value = 3899.8333
$ws.Range("K122").Value = 11699.4999
$ws.Range("M122").Value = -9249.499899999999
